$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 588-611: obstation (588-593) + agripark list/eval/check (594-611) ---
$rowData = @{
  588 = @{ 'C'='序号'; 'D'='index'; 'F'='v99'; 'G'='obstation'; 'H'='list'; 'I'='index'; 'J'='公共网站'; 'K'='国家野外科学观测研究站'; 'L'='名单'; 'M'='序号'; 'O'='v2022.7'; 'P'='科技部网站' }
  589 = @{ 'C'='名称'; 'D'='name'; 'F'='v99'; 'G'='obstation'; 'H'='list'; 'I'='name'; 'J'='公共网站'; 'K'='国家野外科学观测研究站'; 'L'='名单'; 'M'='国家野外站名称'; 'O'='v2022.7'; 'P'='科技部网站' }
  590 = @{ 'C'='依托单位'; 'D'='institution'; 'F'='v99'; 'G'='obstation'; 'H'='list'; 'I'='institution'; 'J'='公共网站'; 'K'='国家野外科学观测研究站'; 'L'='名单'; 'M'='依托单位'; 'O'='v2022.7'; 'P'='科技部网站' }
  591 = @{ 'C'='主管部门'; 'D'='administrator'; 'F'='v99'; 'G'='obstation'; 'H'='list'; 'I'='administrator'; 'J'='公共网站'; 'K'='国家野外科学观测研究站'; 'L'='名单'; 'M'='主管部门'; 'O'='v2022.7'; 'P'='科技部网站' }
  592 = @{ 'C'='公示部门'; 'D'='officer'; 'F'='v99'; 'G'='obstation'; 'H'='list'; 'I'='officer'; 'J'='公共网站'; 'K'='国家野外科学观测研究站'; 'L'='名单'; 'M'='公示部门'; 'O'='v2022.7'; 'P'='科技部网站' }
  593 = @{ 'C'='省份'; 'D'='province'; 'F'='v99'; 'G'='obstation'; 'H'='list'; 'I'='province'; 'J'='公共网站'; 'K'='国家野外科学观测研究站'; 'L'='名单'; 'M'='省份'; 'O'='v2022.7'; 'P'='科技部网站' }
  594 = @{ 'C'='年份'; 'D'='year'; 'F'='v99'; 'G'='agripark'; 'H'='list'; 'I'='year'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='名单'; 'M'='年份'; 'O'='v2022.7'; 'P'='科技部网站' }
  595 = @{ 'C'='序号'; 'D'='index'; 'F'='v99'; 'G'='agripark'; 'H'='list'; 'I'='index'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='名单'; 'M'='序号'; 'O'='v2022.7'; 'P'='科技部网站' }
  596 = @{ 'C'='批次'; 'D'='batch'; 'F'='v99'; 'G'='agripark'; 'H'='list'; 'I'='batch'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='名单'; 'M'='批次'; 'O'='v2022.7'; 'P'='科技部网站' }
  597 = @{ 'C'='名称'; 'D'='name'; 'F'='v99'; 'G'='agripark'; 'H'='list'; 'I'='name'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='名单'; 'M'='名称'; 'O'='v2022.7'; 'P'='科技部网站' }
  598 = @{ 'C'='省份'; 'D'='province'; 'F'='v99'; 'G'='agripark'; 'H'='list'; 'I'='province'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='名单'; 'M'='省份'; 'O'='v2022.7'; 'P'='科技部网站' }
  599 = @{ 'C'='公示部门'; 'D'='officer'; 'F'='v99'; 'G'='agripark'; 'H'='list'; 'I'='officer'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='名单'; 'M'='公示部门'; 'O'='v2022.7'; 'P'='科技部网站' }
  600 = @{ 'C'='年份'; 'D'='year'; 'F'='v99'; 'G'='agripark'; 'H'='eval'; 'I'='year'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='评估'; 'M'='年份'; 'O'='v2022.7'; 'P'='科技部网站' }
  601 = @{ 'C'='序号'; 'D'='index'; 'F'='v99'; 'G'='agripark'; 'H'='eval'; 'I'='index'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='评估'; 'M'='序号'; 'O'='v2022.7'; 'P'='科技部网站' }
  602 = @{ 'C'='名称'; 'D'='name'; 'F'='v99'; 'G'='agripark'; 'H'='eval'; 'I'='name'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='评估'; 'M'='名称'; 'O'='v2022.7'; 'P'='科技部网站' }
  603 = @{ 'C'='名称'; 'D'='name'; 'F'='v99'; 'G'='agripark'; 'H'='eval'; 'I'='result'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='评估'; 'M'='结果'; 'O'='v2022.7'; 'P'='科技部网站' }
  604 = @{ 'C'='省份'; 'D'='province'; 'F'='v99'; 'G'='agripark'; 'H'='eval'; 'I'='province'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='评估'; 'M'='省份'; 'O'='v2022.7'; 'P'='科技部网站' }
  605 = @{ 'C'='公示部门'; 'D'='officer'; 'F'='v99'; 'G'='agripark'; 'H'='eval'; 'I'='officer'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='评估'; 'M'='公示部门'; 'O'='v2022.7'; 'P'='科技部网站' }
  606 = @{ 'C'='年份'; 'D'='year'; 'F'='v99'; 'G'='agripark'; 'H'='check'; 'I'='year'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='验收'; 'M'='年份'; 'O'='v2022.7'; 'P'='科技部网站' }
  607 = @{ 'C'='序号'; 'D'='index'; 'F'='v99'; 'G'='agripark'; 'H'='check'; 'I'='index'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='验收'; 'M'='序号'; 'O'='v2022.7'; 'P'='科技部网站' }
  608 = @{ 'C'='名称'; 'D'='name'; 'F'='v99'; 'G'='agripark'; 'H'='check'; 'I'='name'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='验收'; 'M'='名称'; 'O'='v2022.7'; 'P'='科技部网站' }
  609 = @{ 'C'='名称'; 'D'='name'; 'F'='v99'; 'G'='agripark'; 'H'='check'; 'I'='result'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='验收'; 'M'='结果'; 'O'='v2022.7'; 'P'='科技部网站' }
  610 = @{ 'C'='省份'; 'D'='province'; 'F'='v99'; 'G'='agripark'; 'H'='check'; 'I'='province'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='验收'; 'M'='省份'; 'O'='v2022.7'; 'P'='科技部网站' }
  611 = @{ 'C'='公示部门'; 'D'='officer'; 'F'='v99'; 'G'='agripark'; 'H'='check'; 'I'='officer'; 'J'='公共网站'; 'K'='国家农业科技园区'; 'L'='验收'; 'M'='公示部门'; 'O'='v2022.7'; 'P'='科技部网站' }
}

foreach ($r in $rowData.Keys) {
  $cells = $rowData[$r]
  foreach ($col in $cells.Keys) {
    $ws.Range("$col$r").Value = $cells[$col]
  }
  $ws.Range("A$r").Formula = '=F' + $r + '&"_"&G' + $r + '&"_"&H' + $r + '&"_"&I' + $r
}